$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is a special case: both Invalid (G) and Absent (H) are set to 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows where Total Attendance Count (D) and Real (E) are set to 1
$realRows = @(4, 5, 6, 10, 11, 12, 14)
foreach ($r in $realRows) {
    $ws.Range("D$r").Value = 1
    $ws.Range("E$r").Value = 1
}

# Rows where Absent (H) is set to 1
$absentRows = @(7, 8, 9, 13, 15, 16, 17, 18)
foreach ($r in $absentRows) {
    $ws.Range("H$r").Value = 1
}
